# Daily attendance processing - 2026-02-11 17:16:42 UTC
# Fills in the "Recorded By" identifier for sessions that have now been
# recorded, refreshes the derived attendance counts/percentages that move
# as a result, and flips the two still-pending Endocrinology sessions over
# to "Recorded" now that their registers have come in.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Recorded By" (column G) for sessions newly marked as recorded ---
# A leading apostrophe forces the purely-numeric recorder id to be stored
# as text (matching the existing "Recorded By" column contents).
$recordedByRows = 2,3,4,5,6,7,8,9,10,11,42,43,44,45,46,47,48,49,50,51,87,127,162,163,164,165,167,168,169,170,171,202,203,204,205,207,208,209,210,211,244,245,246,247,248,249,250,251,252,253,284,285,286,287,288,289,290,291,292,293
foreach ($r in $recordedByRows) {
    $ws.Range("G$r").Value = "'776626600547"
}

# --- Class Statistics block (K6:L10) ---
$ws.Range("L6").Value = 279
$ws.Range("L7").Value = 43
$ws.Range("L9").Value = "86.6%"
$ws.Range("L10").Value = "71.4%"

# --- Per-subject summary rows (GASTROENTEROLOGY, B2-A1) ---
$ws.Range("O17").Value = 38
$ws.Range("P17").Value = 2
$ws.Range("R17").Value = "95.0%"
$ws.Range("S17").Value = "70.5%"

$ws.Range("O18").Value = 34
$ws.Range("P18").Value = 6
$ws.Range("R18").Value = "85.0%"
$ws.Range("S18").Value = "70.4%"

# --- Rows 91 & 131: ENDOCRINOLOGY session 10 (11/02/2026) registers came in ---
# Re-style from the "Not Recorded" (pink) row look to the normal row look,
# then update the attendance count and status text.
$ws.Range("A90:I90").Copy()
$ws.Range("A91:I91").PasteSpecial(-4122)
$ws.Range("H91").Value = "16/40"
$ws.Range("I91").Value = "Recorded"

$ws.Range("A130:I130").Copy()
$ws.Range("A131:I131").PasteSpecial(-4122)
$ws.Range("H131").Value = "26/35"
$ws.Range("I131").Value = "Recorded"

$excel.CutCopyMode = $false
